$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (permutation of data rows 2..24)
# i.e. after the edit, row $dest should contain the D/I/J/K/L/M/O/P values
# that row $src currently holds (before the edit).
$map = @{
    2  = 10
    3  = 17
    4  = 20
    5  = 3
    6  = 11
    7  = 12
    8  = 22
    9  = 13
    10 = 8
    11 = 21
    12 = 23
    13 = 14
    14 = 24
    15 = 16
    16 = 6
    17 = 7
    18 = 9
    19 = 18
    20 = 4
    21 = 2
    22 = 5
    23 = 19
    24 = 15
}

# Columns (by index) that vary row to row and need to move along with the
# permutation: D=4, I=9, J=10, K=11, L=12, M=13, O=15, P=16
$cols = 4,9,10,11,12,13,15,16

# First, snapshot all the source values (before any writes happen), since
# the mapping is a permutation and writes would otherwise clobber values
# that still need to be read for a later destination row.
$snapshot = @{}
foreach ($row in $map.Keys) {
    $src = $map[$row]
    if (-not $snapshot.ContainsKey($src)) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Cells.Item($src, $c).Value2
        }
        $snapshot[$src] = $rowVals
    }
}

# Now apply the snapshot values to their destination rows.
foreach ($row in $map.Keys) {
    $src = $map[$row]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($row, $c).Value = $rowVals[$c]
    }
}
